$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel serial date number for 2020-11-30
$rowDate = 44165

# Row 7
$ws.Range("B7").Value = $rowDate
$ws.Range("C7").Value = "Fertigstellung Basis Login System"
$ws.Range("D7").Value = "8:00 Uhr"
$ws.Range("E7").Value = "12:30 Uhr"
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = "Pause"
$ws.Range("H7").Value = 0.5

# Row 8
$ws.Range("B8").Value = $rowDate
$ws.Range("C8").Value = "Meeting"
$ws.Range("D8").Value = "13:00 Uhr"
$ws.Range("E8").Value = "14:00 Uhr"
$ws.Range("F8").Value = 1

# Match the date formatting/style already used by the other date cells
$ws.Range("B6").Copy()
$ws.Range("B7:B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("G8").Select()
